$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 (pushes existing rows 5..21 down to 6..22),
# before touching D4's style so the new row 5 inherits the un-edited formats.
$ws.Rows.Item(5).Insert()

# New row 5 height (task entry row).
$ws.Rows.Item(5).RowHeight = 63.75

# --- Update existing cells with new content ---
$ws.Range("E3").Value = "完成，可以正常usb和网络通信"
$ws.Range("E3").VerticalAlignment = -4108
$ws.Range("E3").WrapText = $true

$ws.Range("D4").Value = "本周完成通讯协议的设计与实现（3.24）"
$ws.Range("D4").VerticalAlignment = -4108
$ws.Range("D4").WrapText = $true

# --- Fill the newly inserted row 5 ---
$ws.Range("A5").Value = 3

$ws.Range("B5").Value = "功能实现：" + [char]10
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("B5").WrapText = $true

$ws.Range("C5").Value = "1、usb 速度测试" + [char]10 + "2、LCD屏显示"
$ws.Range("C5").VerticalAlignment = -4108
$ws.Range("C5").WrapText = $true

# --- Renumber the task index column for the rows pushed down by the insert ---
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11

# --- Update selection ---
$ws.Range("E6").Select() | Out-Null
